# Doing Updates for Financials
# Insert a new first data column (column D) in front of the existing
# yearly columns, shifting the old D:K data right to E:L, and fill the
# new column D with the newest fiscal year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before D; existing D:K data shifts to E:L.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D inherits formatting from column C by
#    default. Copy number-format/style from column E (the old column D,
#    now shifted right) onto the new column D so every row matches the
#    style used by the rest of that row's year columns. Only do this for
#    the rows that actually hold year-column data (skip section-heading
#    rows that have no D:K cells at all, so we don't invent empty cells
#    there).
$dataRowRanges = @("7:35", "38:77", "80:102")
foreach ($rng in $dataRowRanges) {
    $parts = $rng.Split(":")
    $r1 = $parts[0]
    $r2 = $parts[1]
    $ws.Range("E${r1}:E${r2}").Copy()
    $ws.Range("D${r1}:D${r2}").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 3) Populate the new column D with the newest period's values. Row
#    numbers map to the same rows used by the other year columns.
$newColumnD = @{
    7 = 43465
    8 = 711700
    9 = 628400
    10 = 83300
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 670000
    18 = 41700
    20 = 0
    21 = 49400
    22 = 0
    23 = 41700
    24 = 8400
    25 = 0
    26 = 33300
    27 = 33300
    28 = 0
    29 = 500
    30 = 0
    31 = 0
    32 = 0
    33 = 33700
    34 = 0
    35 = 33700
    38 = 43465
    41 = 27000
    42 = 0
    43 = 149100
    44 = 93800
    45 = 3300
    46 = 273200
    47 = 0
    48 = 82800
    49 = 11600
    50 = 0
    51 = 0
    52 = 500
    53 = 0
    54 = 368200
    57 = 98200
    58 = 300
    59 = 24900
    60 = 123400
    61 = 15500
    62 = 1700
    63 = 0
    64 = 0
    65 = 0
    66 = 140600
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 81400
    73 = 0
    74 = 0
    75 = 0
    76 = 227600
    77 = 0
    80 = 43465
    81 = 33700
    83 = 7700
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 21900
    91 = -13300
    92 = 0
    93 = 0
    94 = -13200
    96 = -8200
    97 = 0
    98 = 0
    99 = 0
    100 = -3000
    101 = -600
    102 = 5100
}

foreach ($row in $newColumnD.Keys) {
    $ws.Cells.Item($row, 4).Value = $newColumnD[$row]
}
